$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.74"
$ws.Range("E2").Value = "'-3.72%"
$ws.Range("G2").Value = "'3"
$ws.Range("D3").Value = "'48.60"
$ws.Range("E3").Value = "'3.15%"
$ws.Range("G3").Value = "'3"
$ws.Range("D4").Value = "'5.132"
$ws.Range("E4").Value = "'-2.75%"
$ws.Range("G4").Value = "'3"
$ws.Range("D5").Value = "'0.07746"
$ws.Range("E5").Value = "'-4.44%"
$ws.Range("G5").Value = "'3"
$ws.Range("D6").Value = "'4.455"
$ws.Range("G6").Value = "'3"
$ws.Range("D7").Value = "'1.311"
$ws.Range("E7").Value = "'20.00%"
$ws.Range("G7").Value = "'3"
$ws.Range("D8").Value = "'1.579"
$ws.Range("E8").Value = "'-5.14%"
$ws.Range("G8").Value = "'3"
$ws.Range("D9").Value = "'0.1216"
$ws.Range("E9").Value = "'-7.86%"
$ws.Range("G9").Value = "'3"
$ws.Range("D10").Value = "'0.1933"
$ws.Range("E10").Value = "'-1.43%"
$ws.Range("G10").Value = "'3"
$ws.Range("D11").Value = "'0.04696"
$ws.Range("E11").Value = "'2.70%"
$ws.Range("G11").Value = "'3"
$ws.Range("D12").Value = "'0.09338"
$ws.Range("E12").Value = "'-3.14%"
$ws.Range("G12").Value = "'3"
$ws.Range("E13").Value = "'0.06%"
$ws.Range("G13").Value = "'3"
$ws.Range("D14").Value = "'0.001278"
$ws.Range("E14").Value = "'-3.72%"
$ws.Range("G14").Value = "'3"
$ws.Range("D15").Value = "'0.04170"
$ws.Range("E15").Value = "'-3.05%"
$ws.Range("G15").Value = "'3"
$ws.Range("D16").Value = "'0.005852"
$ws.Range("E16").Value = "'1.39%"
$ws.Range("G16").Value = "'3"
$ws.Range("D17").Value = "'3.346"
$ws.Range("E17").Value = "'-1.17%"
$ws.Range("G17").Value = "'3"
$ws.Range("D18").Value = "'2.281"
$ws.Range("E18").Value = "'-6.25%"
$ws.Range("G18").Value = "'3"
$ws.Range("E19").Value = "'2.94%"
$ws.Range("G19").Value = "'3"
$ws.Range("D20").Value = "'8.115"
$ws.Range("E20").Value = "'-0.54%"
$ws.Range("G20").Value = "'3"
$ws.Range("D21").Value = "'0.1356"
$ws.Range("E21").Value = "'-2.26%"
$ws.Range("G21").Value = "'3"
$ws.Range("D22").Value = "'0.3013"
$ws.Range("E22").Value = "'-1.99%"
$ws.Range("G22").Value = "'3"
$ws.Range("D23").Value = "'0.001268"
$ws.Range("E23").Value = "'-2.83%"
$ws.Range("G23").Value = "'3"
$ws.Range("D24").Value = "'0.004082"
$ws.Range("E24").Value = "'-4.15%"
$ws.Range("G24").Value = "'3"
$ws.Range("D25").Value = "'0.0001351"
$ws.Range("E25").Value = "'0.23%"
$ws.Range("G25").Value = "'3"
$ws.Range("E26").Value = "'-3.85%"
$ws.Range("G26").Value = "'3"
$ws.Range("G27").Value = "'3"
$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("D38").Value = "'0.02565"
$ws.Range("E38").Value = "'-6.89%"
$ws.Range("G38").Value = "'3"
$ws.Range("D39").Value = "'0.05808"
$ws.Range("E39").Value = "'4.81%"
$ws.Range("G39").Value = "'3"
$ws.Range("D40").Value = "'0.01077"
$ws.Range("E40").Value = "'85.97%"
$ws.Range("G40").Value = "'3"
$ws.Range("D41").Value = "'0.007914"
$ws.Range("E41").Value = "'1.75%"
$ws.Range("G41").Value = "'3"
$ws.Range("D42").Value = "'0.1418"
$ws.Range("E42").Value = "'-1.96%"
$ws.Range("G42").Value = "'3"
$ws.Range("D43").Value = "'0.008430"
$ws.Range("E43").Value = "'10.02%"
$ws.Range("G43").Value = "'3"
$ws.Range("D44").Value = "'0.007637"
$ws.Range("E44").Value = "'-13.68%"
$ws.Range("G44").Value = "'3"
$ws.Range("D45").Value = "'0.3355"
$ws.Range("E45").Value = "'-4.28%"
$ws.Range("G45").Value = "'3"
$ws.Range("D46").Value = "'0.00006849"
$ws.Range("E46").Value = "'0.16%"
$ws.Range("G46").Value = "'3"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.24%"
$ws.Range("G47").Value = "'3"
$ws.Range("D48").Value = "'0.05664"
$ws.Range("E48").Value = "'6.26%"
$ws.Range("G48").Value = "'3"
$ws.Range("G49").Value = "'3"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.24%"
$ws.Range("G50").Value = "'3"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.24%"
$ws.Range("G51").Value = "'3"
